# Add a new "Software Licenses" worksheet (with a backing table) to track
# software license assignments, placed after the last existing sheet
# ("Account Information") and made the active/selected sheet.

$wb = $excel.ActiveWorkbook

# New sheet goes after the current last worksheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Software Licenses"

# Header row.
$headers = @("Company", "Software", "Version", "License Key", "Assigned To", "Assigned Date")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Column widths to roughly match the final layout.
$ws.Columns.Item(1).ColumnWidth = 9.6
$ws.Columns.Item(2).ColumnWidth = 9.25
$ws.Columns.Item(3).ColumnWidth = 8.1
$ws.Columns.Item(4).ColumnWidth = 11.25
$ws.Columns.Item(5).ColumnWidth = 11.75
$ws.Columns.Item(6).ColumnWidth = 13.67

# Turn the header row (plus one empty data row) into a real Excel table.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:F2"), $null, 1)
$tbl.Name = "Table8"
$tbl.TableStyle = "TableStyleMedium2"

# Leave the selection on the first data row, matching a freshly created table.
[void]$ws.Range("A2").Select()
